# Updates "想去人数" (F column) and one "最低票价" (G column) values across the
# four worksheets of the 上海-漫展信息 workbook, matching the data refresh
# recorded in commit "Update gh-pages to output generated at 456a3b4".

$wb = $excel.ActiveWorkbook

# --- Sheet 1: 展览 ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("F4").Value = 1204
$ws1.Range("F5").Value = 1648
$ws1.Range("F6").Value = 871
$ws1.Range("F8").Value = 2160
$ws1.Range("F9").Value = 657
$ws1.Range("F10").Value = 537
$ws1.Range("F13").Value = 284
$ws1.Range("F18").Value = 2574
$ws1.Range("F21").Value = 311
$ws1.Range("F22").Value = 1694
$ws1.Range("F25").Value = 538
$ws1.Range("F27").Value = 4455

# --- Sheet 2: 演出 ---
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("F2").Value = 402
$ws2.Range("F3").Value = 377
$ws2.Range("F14").Value = 306
$ws2.Range("F19").Value = 267
$ws2.Range("F25").Value = 228
$ws2.Range("G26").Value = 180

# --- Sheet 3: 本地生活 ---
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("F4").Value = 1394
$ws3.Range("F5").Value = 1769
$ws3.Range("F7").Value = 469
$ws3.Range("F8").Value = 91

# --- Sheet 4: 全部类型 ---
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("F3").Value = 1394
$ws4.Range("F4").Value = 1769
$ws4.Range("F5").Value = 469
$ws4.Range("F6").Value = 377
$ws4.Range("F11").Value = 1204
$ws4.Range("F12").Value = 1648
$ws4.Range("F16").Value = 871
$ws4.Range("F18").Value = 2160
$ws4.Range("F19").Value = 657
$ws4.Range("F20").Value = 537
$ws4.Range("F23").Value = 284
$ws4.Range("F27").Value = 306
$ws4.Range("F33").Value = 150
$ws4.Range("F34").Value = 2574
$ws4.Range("F35").Value = 267
$ws4.Range("F39").Value = 91
$ws4.Range("F41").Value = 1694
$ws4.Range("F44").Value = 538
$ws4.Range("F47").Value = 4455
